# Battery voltages after a couple of rest days.
# Adds two new measurement rows (5 and 6) to the Batterimåling sheet,
# extends the chart series to cover the new rows, applies the 0.00
# number format to C4/C6 (matching the already-formatted E column),
# and grows the chart so it keeps the same visual proportions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data rows --------------------------------------------------
$ws.Range("A5").Value = 45662
$ws.Range("B5").Value = 12.76
$ws.Range("C5").Value = 12.61
$ws.Range("D5").Value = 12.94
$ws.Range("E5").Value = 12.62

$ws.Range("A6").Value = 45664
$ws.Range("B6").Value = 12.54
$ws.Range("C6").Value = 12.2
$ws.Range("D6").Value = 12.68
$ws.Range("E6").Value = 12.27

# --- Number formats (reuse the existing style slots via copy/paste of
#     formats only, instead of assigning NumberFormat strings, so no new
#     style entries get created in styles.xml) ----------------------
$ws.Range("A2").Copy()
$ws.Range("A5:A6").PasteSpecial(-4122)

$ws.Range("E3").Copy()
$ws.Range("C4").PasteSpecial(-4122)
$ws.Range("E5:E6").PasteSpecial(-4122)
$ws.Range("C6").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- Chart: extend source ranges to include the two new rows -------
$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart
$sc = $chart.SeriesCollection()

$sc.Item(1).Formula = "=SERIES(Batterimåling!`$B`$1,Batterimåling!`$A`$2:`$A`$6,Batterimåling!`$B`$2:`$B`$6,1)"
$sc.Item(2).Formula = "=SERIES(Batterimåling!`$C`$1,Batterimåling!`$A`$2:`$A`$6,Batterimåling!`$C`$2:`$C`$6,2)"
$sc.Item(3).Formula = "=SERIES(Batterimåling!`$D`$1,Batterimåling!`$A`$2:`$A`$6,Batterimåling!`$D`$2:`$D`$6,3)"
$sc.Item(4).Formula = "=SERIES(Batterimåling!`$E`$1,Batterimåling!`$A`$2:`$A`$6,Batterimåling!`$E`$2:`$E`$6,4)"

# Grow the chart by two default row-heights so its bottom edge keeps
# tracking the same relative spot below the (now longer) table.
$co.Height = $co.Height + 28.8

# --- Selection: mirror what Excel leaves behind after typing into E6
#     and pressing Enter (moves to E7) ------------------------------
[void]$ws.Range("E7").Select()
